$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.54"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.120"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.13%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07805"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.38%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.267"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.52%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.890"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.24%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.991"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "8.73%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9247"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.31%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1085"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-10.89%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1894"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.09%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08900"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.61%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03326"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.97%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09583"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.98%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001379"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.25%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005955"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.38%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.497"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.75%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.389"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.31%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3417"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.313"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "20.07%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1277"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.64%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.33%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04336"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.76%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.44%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004247"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02152"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.88%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05017"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.06%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007476"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.80%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.19%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008657"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.72%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002115"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008805"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.01%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006544"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.96%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002853"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.85%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.54%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.08%"


Write-Host "Applied cryptos symbol list update."
